$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "961÷9=" "533÷4="
Replace-Text "608÷3=" "527÷5="
Replace-Text "543÷8=" "843÷2="
Replace-Text "604÷9=" "750÷7="
Replace-Text "993÷3=" "315÷5="
Replace-Text "982÷6=" "677÷9="
Replace-Text "399÷6=" "844÷7="
Replace-Text "858÷5=" "822÷7="
Replace-Text "274÷9=" "949÷6="
Replace-Text "928÷9=" "280÷6="
Replace-Text "428÷3=" "816÷4="
Replace-Text "126÷6=" "161÷7="
Replace-Text "656÷8=" "545÷6="
Replace-Text "634÷9=" "150÷9="
Replace-Text "838÷7=" "314÷8="
Replace-Text "942÷3=" "996÷3="
Replace-Text "796÷3=" "425÷9="
Replace-Text "892÷9=" "630÷2="
Replace-Text "875÷9=" "663÷4="
Replace-Text "961÷2=" "743÷5="
Replace-Text "581÷9=" "776÷3="
Replace-Text "106÷6=" "406÷8="
Replace-Text "900÷3=" "850÷4="
Replace-Text "747÷6=" "607÷3="
Replace-Text "205÷7=" "199÷2="
